$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.182.23"
$ws.Range("E2").Value = "  -1.86%  "

$ws.Range("D3").Value = "3.951.72"
$ws.Range("E3").Value = "  -2.75%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'536.84"
$ws.Range("E5").Value = "  +3.03%  "

$ws.Range("D6").Value = "'148.36"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").Value = "3.945.68"
$ws.Range("E7").Value = "  -2.68%  "

$ws.Range("D8").Value = "'0.686"
$ws.Range("E8").Value = "  -5.79%  "

$ws.Range("E10").Value = "  -5.20%  "

$ws.Range("E11").Value = "  -6.22%  "

$ws.Range("D12").Value = "'55.14"
$ws.Range("E12").Value = "  +12.99%  "

$ws.Range("D13").Value = "'0.0000319"
$ws.Range("E13").Value = "  -4.02%  "

$ws.Range("D14").Value = "'10.65"
$ws.Range("E14").Value = "  -4.18%  "

$ws.Range("D15").Value = "4.585.68"
$ws.Range("E15").Value = "  -2.51%  "

$ws.Range("D16").Value = "3.951.72"
$ws.Range("E16").Value = "  -2.97%  "

$ws.Range("D17").Value = "'13.99"
$ws.Range("E17").Value = "  -3.56%  "

$ws.Range("D18").Value = "'20.53"
$ws.Range("E18").Value = "  -4.17%  "

$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("E20").Value = "  -6.09%  "

$ws.Range("D21").Value = "71.109.89"
$ws.Range("E21").Value = "  -1.93%  "

$ws.Range("D22").Value = "'425.29"
$ws.Range("E22").Value = "  -5.07%  "

$ws.Range("D23").Value = "'97.52"
$ws.Range("E23").Value = "  -6.75%  "

$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("E25").Value = "  +5.41%  "

$ws.Range("D26").Value = "'14.62"
$ws.Range("E26").Value = "  -3.24%  "

$ws.Range("D27").Value = "'11.28"
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").Value = "'3.88"
$ws.Range("E28").Value = "  +17.62%  "

$ws.Range("D29").Value = "'10.76"
$ws.Range("E29").Value = "  -3.31%  "

$ws.Range("D30").Value = "'5.90"
$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("D31").Value = "'36.50"
$ws.Range("E31").Value = "  -4.28%  "

$ws.Range("D32").Value = "'7.78"
$ws.Range("E32").Value = "  +17.55%  "

$ws.Range("D33").Value = "'51.09"
$ws.Range("E33").Value = "  +20.95%  "

$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("D35").Value = "'13.34"
$ws.Range("E35").Value = "  -2.73%  "

$ws.Range("D36").Value = "'676.99"
$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("D37").Value = "'65.51"
$ws.Range("E37").Value = "  -3.71%  "

$ws.Range("D38").Value = "'0.442"
$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("E39").Value = "  -5.84%  "

$ws.Range("E41").Value = "  -2.52%  "

$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").Value = "'0.0484"
$ws.Range("E44").Value = "  -3.30%  "

$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").Value = "'10.21"
$ws.Range("E46").Value = "  +3.60%  "

$ws.Range("E47").Value = "  -6.29%  "

$ws.Range("D48").Value = "'2.65"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("E49").Value = "  -3.23%  "

$ws.Range("D50").Value = "'3.00"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("D51").Value = "'145.21"
$ws.Range("E51").Value = "  +1.04%  "
